$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the statistics for 2014/12-2018/12 (rows 2-6) with corrected values,
# and clear the erroneous 2019(E)-2021(E) rows (7-9) which should have no data.

# Row 2
$ws.Range("D2").Value = 1006
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 115
$ws.Range("H2").Value = 97
$ws.Range("I2").Value = 103
$ws.Range("J2").Value = -6
$ws.Range("K2").Value = 2527
$ws.Range("L2").Value = 684
$ws.Range("M2").Value = 1843
$ws.Range("N2").Value = 1827
$ws.Range("O2").Value = 16
$ws.Range("P2").Value = 41
$ws.Range("Q2").Value = 53
$ws.Range("R2").Value = -96
$ws.Range("S2").Value = -30
$ws.Range("T2").Value = 28
$ws.Range("U2").Value = 25
$ws.Range("V2").Value = 125
$ws.Range("W2").Value = 0.8
$ws.Range("X2").Value = 9.609999999999999
$ws.Range("Y2").Value = 5.75
$ws.Range("Z2").Value = 3.84
$ws.Range("AA2").Value = 37.15
$ws.Range("AB2").Value = 4360.47
$ws.Range("AC2").Value = 1269
$ws.Range("AD2").Value = 9.300000000000001
$ws.Range("AE2").Value = 22779
$ws.Range("AF2").Value = 0.52
$ws.Range("AG2").Value = 80
$ws.Range("AH2").Value = 0.68
$ws.Range("AI2").Value = 6.22
$ws.Range("AJ2").Value = 8126314

# Row 3
$ws.Range("D3").Value = 924
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = 31
$ws.Range("G3").Value = 75
$ws.Range("H3").Value = 84
$ws.Range("I3").Value = 86
$ws.Range("J3").Value = -2
$ws.Range("K3").Value = 2441
$ws.Range("L3").Value = 519
$ws.Range("M3").Value = 1922
$ws.Range("N3").Value = 1909
$ws.Range("O3").Value = 13
$ws.Range("P3").Value = 41
$ws.Range("Q3").Value = 86
$ws.Range("R3").Value = -50
$ws.Range("S3").Value = -12
$ws.Range("T3").Value = 25
$ws.Range("U3").Value = 60
$ws.Range("V3").Value = 119
$ws.Range("W3").Value = 3.38
$ws.Range("X3").Value = 9.09
$ws.Range("Y3").Value = 4.62
$ws.Range("Z3").Value = 3.38
$ws.Range("AA3").Value = 27.02
$ws.Range("AB3").Value = 4570.29
$ws.Range("AC3").Value = 1061
$ws.Range("AD3").Value = 9.85
$ws.Range("AE3").Value = 23804
$ws.Range("AF3").Value = 0.44
$ws.Range("AG3").Value = 100
$ws.Range("AH3").Value = 0.96
$ws.Range("AI3").Value = 9.300000000000001
$ws.Range("AJ3").Value = 8126314

# Row 4
$ws.Range("D4").Value = 847
$ws.Range("E4").Value = 19
$ws.Range("F4").Value = 19
$ws.Range("G4").Value = 108
$ws.Range("H4").Value = 78
$ws.Range("I4").Value = 77
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2527
$ws.Range("L4").Value = 540
$ws.Range("M4").Value = 1987
$ws.Range("N4").Value = 1973
$ws.Range("O4").Value = 14
$ws.Range("P4").Value = 41
$ws.Range("Q4").Value = 92
$ws.Range("R4").Value = -6
$ws.Range("S4").Value = -15
$ws.Range("T4").Value = 9
$ws.Range("U4").Value = 84
$ws.Range("V4").Value = 112
$ws.Range("W4").Value = 2.2
$ws.Range("X4").Value = 9.15
$ws.Range("Y4").Value = 3.97
$ws.Range("Z4").Value = 3.12
$ws.Range("AA4").Value = 27.17
$ws.Range("AB4").Value = 4715.66
$ws.Range("AC4").Value = 948
$ws.Range("AD4").Value = 11.97
$ws.Range("AE4").Value = 24604
$ws.Range("AF4").Value = 0.46
$ws.Range("AG4").Value = 90
$ws.Range("AH4").Value = 0.79
$ws.Range("AI4").Value = 9.369999999999999
$ws.Range("AJ4").Value = 8126314

# Row 5
$ws.Range("D5").Value = 918
$ws.Range("E5").Value = 43
$ws.Range("F5").Value = 43
$ws.Range("G5").Value = 140
$ws.Range("H5").Value = 110
$ws.Range("I5").Value = 112
$ws.Range("J5").Value = -2
$ws.Range("K5").Value = 2631
$ws.Range("L5").Value = 549
$ws.Range("M5").Value = 2082
$ws.Range("N5").Value = 2069
$ws.Range("O5").Value = 13
$ws.Range("P5").Value = 41
$ws.Range("Q5").Value = 151
$ws.Range("R5").Value = -216
$ws.Range("S5").Value = -22
$ws.Range("T5").Value = 50
$ws.Range("U5").Value = 101
$ws.Range("V5").Value = 94
$ws.Range("W5").Value = 4.65
$ws.Range("X5").Value = 11.97
$ws.Range("Y5").Value = 5.55
$ws.Range("Z5").Value = 4.26
$ws.Range("AA5").Value = 26.36
$ws.Range("AB5").Value = 4969.44
$ws.Range("AC5").Value = 1379
$ws.Range("AD5").Value = 8.300000000000001
$ws.Range("AE5").Value = 25805
$ws.Range("AF5").Value = 0.44
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 0.87
$ws.Range("AI5").Value = 7.15
$ws.Range("AJ5").Value = 8126314

# Row 6
$ws.Range("D6").Value = 809
$ws.Range("E6").Value = -11
$ws.Range("F6").Value = -11
$ws.Range("G6").Value = 80
$ws.Range("H6").Value = 63
$ws.Range("I6").Value = 63
$ws.Range("K6").Value = 2602
$ws.Range("L6").Value = 463
$ws.Range("M6").Value = 2139
$ws.Range("N6").Value = 2125
$ws.Range("P6").Value = 41
$ws.Range("Q6").Value = 15
$ws.Range("R6").Value = 54
$ws.Range("S6").Value = -11
$ws.Range("T6").Value = 37
$ws.Range("U6").Value = -22
$ws.Range("V6").Value = 91
$ws.Range("W6").Value = -1.4
$ws.Range("X6").Value = 7.83
$ws.Range("Y6").Value = 2.99
$ws.Range("Z6").Value = 2.42
$ws.Range("AA6").Value = 21.64
$ws.Range("AB6").Value = 5088.52
$ws.Range("AC6").Value = 771
$ws.Range("AD6").Value = 11.81
$ws.Range("AE6").Value = 26505
$ws.Range("AF6").Value = 0.34
$ws.Range("AG6").Value = 80
$ws.Range("AH6").Value = 0.88
$ws.Range("AI6").Value = 10.24
$ws.Range("AJ6").Value = 8126314

# Rows 7-9: clear all data cells (D:AJ) - these projected years are removed
$ws.Range("D7:AJ9").ClearContents()
